{"js": "// Remove the stray italic \"2 Samuel\" paragraph that sits right after the\n// \"2SA\" (Heading 2) short-code heading. The whole paragraph (its run(s)\n// and its paragraph mark) is deleted; the \"2SA\" heading paragraph and the\n// following blank/space paragraph are left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \"2 Samuel\") {\n    p.font.load(\"italic\");\n    target = target || [];\n    target.push(p);\n  }\n}\n\nif (target && target.length) {\n  await context.sync();\n  for (const p of target) {\n    if (p.font.italic === true) {\n      p.delete();\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the stray italic \"2 Samuel\" paragraph that sits right after the\n# \"2SA\" (Heading 2) short-code heading. The whole paragraph (its run(s)\n# and its paragraph mark) is deleted; the \"2SA\" heading paragraph and the\n# following blank/space paragraph are left untouched.\n\n$d = $word.ActiveDocument\n\n$targets = New-Object System.Collections.ArrayList\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $text = $r.Text.Trim()\n    if ($r.Italic -and ($text -eq \"2 Samuel\")) {\n        [void]$targets.Add($p)\n    }\n}\n\nfor ($i = $targets.Count - 1; $i -ge 0; $i--) {\n    $targets[$i].Range.Delete()\n}\n"}
